$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")
$ws.Activate()

$ws.Range("A3").Value = 0.02
$ws.Range("B3").Value = 0.134
$ws.Range("A4").Value = 0.04
$ws.Range("B4").Value = 0.23
$ws.Range("A5").Value = 0.06
$ws.Range("B5").Value = 0.304
$ws.Range("A6").Value = 0.08
$ws.Range("B6").Value = 0.365
$ws.Range("A7").Value = 0.1
$ws.Range("B7").Value = 0.418
$ws.Range("A8").Value = 0.15
$ws.Range("B8").Value = 0.517
$ws.Range("A9").Value = 0.2
$ws.Range("B9").Value = 0.579
$ws.Range("A10").Value = 0.3
$ws.Range("B10").Value = 0.665
$ws.Range("A11").Value = 0.4
$ws.Range("B11").Value = 0.729
$ws.Range("A12").Value = 0.5
$ws.Range("B12").Value = 0.779
$ws.Range("A13").Value = 0.6
$ws.Range("B13").Value = 0.825
$ws.Range("A14").Value = 0.7
$ws.Range("B14").Value = 0.87
$ws.Range("A15").Value = 0.8
$ws.Range("B15").Value = 0.915
$ws.Range("A16").Value = 0.9
$ws.Range("B16").Value = 0.958
$ws.Range("A17").Value = 0.95
$ws.Range("B17").Value = 0.979
$ws.Range("A18").Value = 1
$ws.Range("B18").Value = 1

[void]$ws.Range("B18").Select()
